$d = $word.ActiveDocument

# Locate the last paragraph in the document ("(Unless you really want a goat)")
$lastPara = $d.Paragraphs.Last
$countBefore = $d.Paragraphs.Count

# Insert a placeholder empty paragraph right after it; InsertXML replaces the
# *whole* contents of the range it is called on (including the paragraph
# mark), so we give it a throw-away paragraph to consume rather than letting
# it eat the existing "goat" paragraph.
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$placeholder = $d.Paragraphs.Last
$r2 = $placeholder.Range
$r2.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"><w:body><w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr>
          <w:color w:val="009900"/>
          <w:lang w:val="en-IN" w:eastAsia="en-IN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="FF6600"/>
        </w:rPr>
        <w:t>Puzzle to place seven match sticks in order that each match stick crosses the other 6</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="009900"/>
          <w:lang w:val="en-IN" w:eastAsia="en-IN"/>
        </w:rPr>
        <w:pict>
          <v:group id="shape_0" style="position:absolute;margin-left:0pt;margin-top:0pt;width:179.95pt;height:170.95pt" coordorigin="0,0" coordsize="3599,3419">
            <v:rect id="shape_0" stroked="f" style="position:absolute;left:0;top:0;width:3598;height:3418;mso-position-horizontal-relative:char">
              <v:wrap v:type="none"/>
              <v:fill on="false" detectmouseclick="t"/>
              <v:stroke color="#3465a4" joinstyle="round" endcap="flat"/>
            </v:rect>
            <v:line id="shape_0" from="538,540" to="1796,1978" stroked="t" style="position:absolute;flip:x;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="1799,540" to="2877,1978" stroked="t" style="position:absolute;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="538,1979" to="2876,1979" stroked="t" style="position:absolute;flip:x;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="1799,540" to="1799,3418" stroked="t" style="position:absolute;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="539,900" to="3597,1977" stroked="t" style="position:absolute;flip:y;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="718,1079" to="2876,1977" stroked="t" style="position:absolute;flip:xy;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
            <v:line id="shape_0" from="1078,540" to="1796,2518" stroked="t" style="position:absolute;flip:x;mso-position-horizontal-relative:char">
              <v:stroke color="black" weight="9360" joinstyle="miter" endcap="square"/>
              <v:fill on="false" detectmouseclick="t"/>
            </v:line>
          </v:group>
        </w:pict>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr>
          <w:color w:val="FF6600"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="FF6600"/>
        </w:rPr>
        <w:t>Make 1000 using 8 8's</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
        <w:t>There are 2 solutions :</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
        <w:t>&gt; 888 + 88 + 8 + 8 + 8 and</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:color w:val="009900"/>
        </w:rPr>
        <w:t>&gt; ( 8888-888 ) / 8</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Normal"/>
        <w:spacing w:before="0" w:after="0"/>
        <w:rPr/>
      </w:pPr>
      <w:r>
        <w:rPr/>
      </w:r>
    </w:p></w:body></w:document>
'@

$r2.InsertXML($xml)

$countAfter = $d.Paragraphs.Count
$numInserted = $countAfter - $countBefore

# The WordprocessingML <w:spacing w:before="0" w:after="0"/> written in the
# inserted fragment gets normalised away by InsertXML, so reapply it
# explicitly through the paragraph format object model for every paragraph
# that was just inserted.
for ($i = $countBefore + 1; $i -le $countAfter; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.SpaceBefore = 0
    $p.SpaceAfter = 0
}

Write-Host "Paragraphs before:" $countBefore "after:" $countAfter "inserted:" $numInserted
